# Improved accuracy of stimulus presentation time-logging
# Updates task-order worksheet names (embedded timestamps) and the
# CSV-filename / condition-order values they contain.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (new embedded timestamps) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555040398927"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555072088923"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555072148976"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555072738988"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555073508916"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555040108953.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555040238945.csv"
$ws1.Range("B4").Value = "go_stims-16512555040248947.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555040388925.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512555056738968.csv"
$ws2.Range("B3").Value = "OB-16512555055208926.csv"
$ws2.Range("B4").Value = "ZB-match_4-16512555051328924.csv"
$ws2.Range("B5").Value = "TB-16512555060848963.csv"
$ws2.Range("B6").Value = "OB-16512555051628928.csv"
$ws2.Range("B7").Value = "TB-16512555071908922.csv"
$ws2.Range("B8").Value = "TB-16512555068078933.csv"
$ws2.Range("B9").Value = "ZB-match_1-16512555048598917.csv"
$ws2.Range("B10").Value = "ZB-match_1-16512555041478963.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555072398922.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555072178936.csv"
$ws4.Range("B4").Value = "MM_stims-16512555072558925.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555072398922.csv"
$ws4.Range("B6").Value = "MM_stims-16512555072718935.csv"
$ws4.Range("B7").Value = "ZM_stims-1651255507256895.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512555072798939.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555073348927.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555073188946.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555073038914.csv"
